# Apply the "Today's update - 13Dec" data refresh to the DateofDeath sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DateofDeath")

# --- Corrections to previously entered rows -----------------------------
# D267 (Probable Deaths) revised upward
$ws.Range("D267").Value = 1

# B270 / B272 (Confirmed Deaths) revised upward
$ws.Range("B270").Value = 44
$ws.Range("B272").Value = 44

# B275 (Confirmed Deaths) revised upward
$ws.Range("B275").Value = 40

# B276 / D276 revised upward
$ws.Range("B276").Value = 48
$ws.Range("D276").Value = 2

# --- New rows for 2020-12-10 and 2020-12-11 ------------------------------
# Row 277: 2020-12-10
$ws.Range("A277").Value = 44175
$ws.Range("A277").Style = $ws.Range("A276").Style
$ws.Range("B277").Value = 46
$ws.Range("C276").AutoFill($ws.Range("C276:C277"), 0) | Out-Null
$ws.Range("D277").Value = 1
$ws.Range("E276").AutoFill($ws.Range("E276:E277"), 0) | Out-Null
$ws.Range("F276").AutoFill($ws.Range("F276:F277"), 0) | Out-Null

# Row 278: 2020-12-11
$ws.Range("A278").Value = 44176
$ws.Range("A278").Style = $ws.Range("A277").Style
$ws.Range("B278").Value = 24
$ws.Range("C277").AutoFill($ws.Range("C277:C278"), 0) | Out-Null
$ws.Range("D278").Value = 0
$ws.Range("E277").AutoFill($ws.Range("E277:E278"), 0) | Out-Null
$ws.Range("F277").AutoFill($ws.Range("F277:F278"), 0) | Out-Null

# --- Keep frozen pane / selection in sync with the new data extent ------
$ws.Range("D2:D278").Select() | Out-Null

$excel.ActiveWindow.Panes.Item(4).TopLeftCell = $ws.Range("B269")

$wb.Application.CalculateFullRebuild()
